$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 — shifts the existing rows 4..97 down to 5..98
$ws.Rows.Item(4).EntireRow.Insert()

# Populate the newly inserted row 4 with the new weekly record
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 45237
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100103
$ws.Range("H4").Value = "Frutos de hueso (carozo)"
$ws.Range("I4").Value = 100103004
$ws.Range("J4").Value = "Durazno"
$ws.Range("K4").Value = "Florida King"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 550
$ws.Range("N4").Value = 22000
$ws.Range("O4").Value = 33000
$ws.Range("P4").Value = 27045
$ws.Range("Q4").Value = "$/bandeja 18 kilos granel"
$ws.Range("R4").Value = "Región de Coquimbo"
$ws.Range("S4").Value = 1502
$ws.Range("T4").Value = 18
